$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 248 - this shifts rows 248..348 down to 249..349
$ws.Rows.Item(248).Insert()

# Populate the newly inserted row 248 with the new data record.
$ws.Cells.Item(248, 1).Value = 5
$ws.Cells.Item(248, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(248, 3).Value = "Maule"
$ws.Cells.Item(248, 4).Value = 44755
$ws.Cells.Item(248, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(248, 5).Value = 7
$ws.Cells.Item(248, 6).Value = 100114014
$ws.Cells.Item(248, 7).Value = "Betarraga"
$ws.Cells.Item(248, 8).Value = "Sin especificar"
$ws.Cells.Item(248, 9).Value = "Primera"
$ws.Cells.Item(248, 10).Value = 3000
$ws.Cells.Item(248, 11).Value = 750
$ws.Cells.Item(248, 12).Value = 750
$ws.Cells.Item(248, 13).Value = 750
$ws.Cells.Item(248, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(248, 15).Value = "Región del Maule"
$ws.Cells.Item(248, 16).Value = 150
$ws.Cells.Item(248, 17).Value = 5
$ws.Cells.Item(248, 18).Value = "Hortaliza"
